$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "rm 2 US & add Native Son"
# The US (row 269) entry "The Color Purple" / "Alice Walker" becomes
# "Native Son" / "Richard Wright", and is now marked as read.
$ws.Range("C269").Value = "Native Son"
$ws.Range("D269").Value = "Richard Wright"
$ws.Range("E269").Value = 1

# Remove the two US rows for "Catcher in the Rye" (J. D. Salinger) and
# "Catch 22" (Joseph Heller). Row 270 is "Catcher in the Rye"; once it is
# removed, "Catch 22" (originally row 272) shifts up to row 271.
$ws.Rows.Item(270).Delete()
$ws.Rows.Item(271).Delete()
